$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.440.65'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '3.691.55'
$ws.Range('E3').Value = '  -2.71%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '687.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.63%  '
$ws.Range('D7').Value = '3.690.89'
$ws.Range('E7').Value = '  -2.68%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('E9').Value = '  -5.76%  '
$ws.Range('E10').Value = '  -8.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.20'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.50%  '
$ws.Range('E12').Value = '  -8.96%  '
$ws.Range('E13').Value = '  -6.60%  '
$ws.Range('D14').Value = '4.314.39'
$ws.Range('E14').Value = '  -2.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -10.14%  '
$ws.Range('D16').Value = '3.696.81'
$ws.Range('E16').Value = '  -3.51%  '
$ws.Range('D17').Value = '69.473.98'
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -9.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -10.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '469.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.26%  '
$ws.Range('E23').Value = '  -9.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.77%  '
$ws.Range('D25').Value = '3.838.10'
$ws.Range('E25').Value = '  -2.69%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000126'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -10.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -13.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.57%  '
$ws.Range('E30').Value = '  -8.67%  '
$ws.Range('E31').Value = '  -12.55%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -10.61%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.65'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.80%  '
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.80'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.87%  '
$ws.Range('E36').Value = '  -6.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.21'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -11.90%  '
$ws.Range('E38').Value = '  -7.29%  '
$ws.Range('E39').Value = '  -3.62%  '
$ws.Range('E41').Value = '  -10.01%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '167.34'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.943'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.93'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.21%  '
$ws.Range('E46').Value = '  -13.59%  '
$ws.Range('E47').Value = '  -4.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '28.68'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.88%  '
$ws.Range('E49').Value = '  -3.81%  '
$ws.Range('E50').Value = '  -8.94%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.82'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -9.34%  '
